# Commit: Changed the app name in the xlsx file to match what's in the test cases code
#
# Rename the worksheet from "adactin" to "com.adactin.hotelapp" so it
# matches the name used in the test automation source code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "com.adactin.hotelapp"

# Renaming the sheet updates any formula/reference that pointed to it
# (e.g. the hidden _FilterDatabase defined name used by the AutoFilter),
# but the sheet name now contains a "." so Excel needs to quote it with
# single quotes in the R1C1-less reference. Make sure that is reflected.
foreach ($name in $wb.Names) {
    $name.RefersTo = "='com.adactin.hotelapp'!`$A`$1:`$I`$20"
}
